$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.425.10"
$ws.Cells.Item(2, 5).Value = "'  -0.23%  "
$ws.Cells.Item(3, 4).Value = "'1.805.94"
$ws.Cells.Item(3, 5).Value = "'  +0.10%  "
$ws.Cells.Item(4, 4).Value = "'1.009"
$ws.Cells.Item(5, 5).Value = "'  +0.26%  "
$ws.Cells.Item(6, 4).Value = "'307.09"
$ws.Cells.Item(6, 5).Value = "'  -0.31%  "
$ws.Cells.Item(7, 4).Value = "'0.4523"
$ws.Cells.Item(7, 5).Value = "'  -0.15%  "
$ws.Cells.Item(8, 4).Value = "'0.3598"
$ws.Cells.Item(8, 5).Value = "'  -1.54%  "
$ws.Cells.Item(9, 4).Value = "'46.13"
$ws.Cells.Item(9, 5).Value = "'  +2.12%  "
$ws.Cells.Item(10, 4).Value = "'0.07073"
$ws.Cells.Item(10, 5).Value = "'  -0.37%  "
$ws.Cells.Item(11, 4).Value = "'0.8902"
$ws.Cells.Item(11, 5).Value = "'  +1.97%  "
$ws.Cells.Item(12, 4).Value = "'0.07813"
$ws.Cells.Item(12, 5).Value = "'  +0.36%  "
$ws.Cells.Item(13, 5).Value = "'  +0.80%  "
$ws.Cells.Item(14, 4).Value = "'1.847.14"
$ws.Cells.Item(14, 5).Value = "'  +1.41%  "
$ws.Cells.Item(15, 4).Value = "'5.288"
$ws.Cells.Item(15, 5).Value = "'  +0.37%  "
$ws.Cells.Item(16, 4).Value = "'6.313"
$ws.Cells.Item(16, 5).Value = "'  -0.31%  "
$ws.Cells.Item(17, 4).Value = "'85.38"
$ws.Cells.Item(17, 5).Value = "'  -1.25%  "
$ws.Cells.Item(18, 5).Value = "'  +0.28%  "
$ws.Cells.Item(19, 4).Value = "'0.000008509"
$ws.Cells.Item(19, 5).Value = "'  -0.38%  "
$ws.Cells.Item(20, 5).Value = "'  +0.13%  "
$ws.Cells.Item(21, 4).Value = "'26.460.78"
$ws.Cells.Item(21, 5).Value = "'  -0.26%  "
$ws.Cells.Item(22, 5).Value = "'  -0.08%  "
$ws.Cells.Item(23, 4).Value = "'4.965"
$ws.Cells.Item(23, 5).Value = "'  +0.06%  "
$ws.Cells.Item(24, 4).Value = "'2.036.05"
$ws.Cells.Item(24, 5).Value = "'  -1.33%  "
$ws.Cells.Item(26, 5).Value = "'  -1.08%  "
$ws.Cells.Item(27, 4).Value = "'152.97"
$ws.Cells.Item(27, 5).Value = "'  +1.69%  "
$ws.Cells.Item(28, 5).Value = "'  -0.22%  "
$ws.Cells.Item(29, 4).Value = "'2.071"
$ws.Cells.Item(29, 5).Value = "'  +3.78%  "
$ws.Cells.Item(30, 4).Value = "'112.12"
$ws.Cells.Item(30, 5).Value = "'  -0.93%  "
$ws.Cells.Item(31, 4).Value = "'4.848"
$ws.Cells.Item(31, 5).Value = "'  -0.54%  "
$ws.Cells.Item(32, 4).Value = "'0.08700"
$ws.Cells.Item(32, 5).Value = "'  +0.53%  "
$ws.Cells.Item(33, 5).Value = "'  +0.22%  "
$ws.Cells.Item(34, 4).Value = "'2.813"
$ws.Cells.Item(34, 5).Value = "'  +11.45%  "
$ws.Cells.Item(35, 4).Value = "'4.460"
$ws.Cells.Item(35, 5).Value = "'  +0.53%  "
$ws.Cells.Item(36, 4).Value = "'0.7233"
$ws.Cells.Item(36, 5).Value = "'  -0.70%  "
$ws.Cells.Item(37, 5).Value = "'  -0.82%  "
$ws.Cells.Item(38, 4).Value = "'1.078"
$ws.Cells.Item(38, 5).Value = "'  +0.01%  "
$ws.Cells.Item(39, 4).Value = "'0.01936"
$ws.Cells.Item(40, 4).Value = "'2.918"
$ws.Cells.Item(40, 5).Value = "'  +1.91%  "
$ws.Cells.Item(41, 4).Value = "'0.05119"
$ws.Cells.Item(41, 5).Value = "'  +0.89%  "
$ws.Cells.Item(42, 4).Value = "'0.5080"
$ws.Cells.Item(42, 5).Value = "'  +3.55%  "
$ws.Cells.Item(43, 4).Value = "'6.779"
$ws.Cells.Item(43, 5).Value = "'  -1.75%  "
$ws.Cells.Item(44, 4).Value = "'0.1511"
$ws.Cells.Item(44, 5).Value = "'  -3.64%  "
$ws.Cells.Item(45, 4).Value = "'8.019"
$ws.Cells.Item(45, 5).Value = "'  -1.34%  "
$ws.Cells.Item(46, 4).Value = "'1.009"
$ws.Cells.Item(46, 5).Value = "'  +0.29%  "
$ws.Cells.Item(47, 4).Value = "'0.4672"
$ws.Cells.Item(47, 5).Value = "'  +1.66%  "
$ws.Cells.Item(48, 4).Value = "'9.943"
$ws.Cells.Item(48, 5).Value = "'  +0.25%  "
$ws.Cells.Item(49, 4).Value = "'100.29"
$ws.Cells.Item(49, 5).Value = "'  -1.20%  "
$ws.Cells.Item(50, 4).Value = "'1.578"
$ws.Cells.Item(50, 5).Value = "'  +0.00%  "
$ws.Cells.Item(51, 4).Value = "'0.05980"
$ws.Cells.Item(51, 5).Value = "'  -0.20%  "
